$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D stores prices as plain text (some values use two "." as thousands
# and decimal separators, so the whole column is text, not numbers). Pre-format
# the cells whose new value would otherwise look like a pure number to Excel so
# they keep being stored as text instead of being auto-converted to a number.
$textCells = @("D5", "D7", "D8", "D9", "D10", "D11", "D15", "D17", "D19", "D21", "D22", "D24", "D25", "D26", "D28", "D30", "D31", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D43", "D44", "D45", "D46", "D47", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "28.654.59"
$ws.Range("E2").Value = "  +1.07%  "

$ws.Range("D3").Value = "1.805.56"

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "317.71"
$ws.Range("E5").Value = "  -0.28%  "

$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").Value = "0.5464"
$ws.Range("E7").Value = "  -4.15%  "

$ws.Range("D8").Value = "0.3811"
$ws.Range("E8").Value = "  -1.65%  "

$ws.Range("D9").Value = "0.07521"
$ws.Range("E9").Value = "  -1.15%  "

$ws.Range("D10").Value = "42.42"
$ws.Range("E10").Value = "  -0.86%  "

$ws.Range("D11").Value = "1.116"
$ws.Range("E11").Value = "  -2.02%  "

$ws.Range("E12").Value = "  +0.00%  "

$ws.Range("E13").Value = "  -2.19%  "

$ws.Range("E14").Value = "  -1.51%  "

$ws.Range("D15").Value = "7.407"
$ws.Range("E15").Value = "  +1.85%  "

$ws.Range("D16").Value = "1.795.67"
$ws.Range("E16").Value = "  -0.90%  "

$ws.Range("D17").Value = "90.25"
$ws.Range("E17").Value = "  -1.82%  "

$ws.Range("E18").Value = "  -0.59%  "

$ws.Range("D19").Value = "0.06487"
$ws.Range("E19").Value = "  +0.15%  "

$ws.Range("E20").Value = "  -0.05%  "

$ws.Range("D21").Value = "17.38"
$ws.Range("E21").Value = "  +0.42%  "

$ws.Range("D22").Value = "5.949"
$ws.Range("E22").Value = "  -0.95%  "

$ws.Range("D23").Value = "28.661.34"
$ws.Range("E23").Value = "  +1.04%  "

$ws.Range("D24").Value = "11.13"
$ws.Range("E24").Value = "  -1.72%  "

$ws.Range("D25").Value = "2.093"
$ws.Range("E25").Value = "  -2.24%  "

$ws.Range("D26").Value = "160.76"
$ws.Range("E26").Value = "  +1.57%  "

$ws.Range("E27").Value = "  -1.64%  "

$ws.Range("D28").Value = "2.366"
$ws.Range("E28").Value = "  -3.70%  "

$ws.Range("D29").Value = "1.998.47"
$ws.Range("E29").Value = "  -1.07%  "

$ws.Range("D30").Value = "123.35"
$ws.Range("E30").Value = "  -0.65%  "

$ws.Range("D31").Value = "1.122"
$ws.Range("E31").Value = "  -3.16%  "

$ws.Range("E32").Value = "  -1.08%  "

$ws.Range("D33").Value = "5.648"
$ws.Range("E33").Value = "  -2.53%  "

$ws.Range("D34").Value = "3.683"
$ws.Range("E34").Value = "  +1.43%  "

$ws.Range("D35").Value = "0.06645"
$ws.Range("E35").Value = "  +8.39%  "

$ws.Range("D36").Value = "0.2262"
$ws.Range("E36").Value = "  +2.21%  "

$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").Value = "8.822"
$ws.Range("E37").Value = "  -1.37%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.02304"
$ws.Range("E38").Value = "  -0.71%  "

$ws.Range("D39").Value = "5.032"
$ws.Range("E39").Value = "  -0.24%  "

$ws.Range("D40").Value = "0.6258"
$ws.Range("E40").Value = "  -2.31%  "

$ws.Range("E41").Value = "  -3.50%  "

$ws.Range("E42").Value = "  +2.69%  "

$ws.Range("D43").Value = "1.436"
$ws.Range("E43").Value = "  +4.06%  "

$ws.Range("D44").Value = "13.29"
$ws.Range("E44").Value = "  -1.25%  "

$ws.Range("D45").Value = "0.5867"
$ws.Range("E45").Value = "  -2.35%  "

$ws.Range("D46").Value = "3.697"
$ws.Range("E46").Value = "  -0.04%  "

$ws.Range("D47").Value = "126.91"
$ws.Range("E47").Value = "  +3.52%  "

$ws.Range("D48").Value = "1.951"
$ws.Range("E48").Value = "  +0.08%  "

$ws.Range("E49").Value = "  +0.90%  "

$ws.Range("D50").Value = "0.06898"
$ws.Range("E50").Value = "  +0.34%  "

$ws.Range("D51").Value = "72.34"
$ws.Range("E51").Value = "  -1.14%  "
